$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Drools rule table was migrated from the "ITR" model object onto the
# "RootJosonObject" domain wrapper (with the totalTI path now nested
# under itr3), so update the three cells that reference it.
$ws.Range("B2").Value = "com.openhack.dev.domain.RootJosonObject"
$ws.Range("B7").Value = '$itrObject: RootJosonObject'
$ws.Range("B8").Value = '$itrObject.itr3.partBTI.totalTI > $param'

# Leave the selection on the condition cell, matching the author's save.
$ws.Range("B8").Select()
